$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to all D/E/G cells being modified so Excel keeps them as text
# (matching the original inline-string text cells), then restore default style afterward.
$changeRange = $ws.Range("D2:G51")
$changeRange.NumberFormat = "@"

$ws.Range("D2").Value = "326.64"
$ws.Range("E2").Value = "-1.08%"
$ws.Range("G2").Value = "4"
$ws.Range("D3").Value = "45.03"
$ws.Range("E3").Value = "2.47%"
$ws.Range("G3").Value = "4"
$ws.Range("D4").Value = "5.568"
$ws.Range("E4").Value = "-5.23%"
$ws.Range("G4").Value = "4"
$ws.Range("D5").Value = "0.08084"
$ws.Range("E5").Value = "-2.66%"
$ws.Range("G5").Value = "4"
$ws.Range("D6").Value = "8.715"
$ws.Range("E6").Value = "-0.82%"
$ws.Range("G6").Value = "4"
$ws.Range("D7").Value = "4.333"
$ws.Range("E7").Value = "-3.78%"
$ws.Range("G7").Value = "4"
$ws.Range("D8").Value = "1.895"
$ws.Range("E8").Value = "-2.87%"
$ws.Range("G8").Value = "4"
$ws.Range("D9").Value = "2.721"
$ws.Range("E9").Value = "-6.76%"
$ws.Range("G9").Value = "4"
$ws.Range("D10").Value = "0.9485"
$ws.Range("E10").Value = "1.97%"
$ws.Range("G10").Value = "4"
$ws.Range("D11").Value = "0.1159"
$ws.Range("E11").Value = "-7.49%"
$ws.Range("G11").Value = "4"
$ws.Range("D12").Value = "0.1896"
$ws.Range("E12").Value = "-2.33%"
$ws.Range("G12").Value = "4"
$ws.Range("D13").Value = "0.1018"
$ws.Range("E13").Value = "7.75%"
$ws.Range("G13").Value = "4"
$ws.Range("E14").Value = "5.23%"
$ws.Range("G14").Value = "4"
$ws.Range("D15").Value = "0.1063"
$ws.Range("E15").Value = "-0.05%"
$ws.Range("G15").Value = "4"
$ws.Range("D16").Value = "0.001290"
$ws.Range("E16").Value = "-1.06%"
$ws.Range("G16").Value = "4"
$ws.Range("D17").Value = "0.005948"
$ws.Range("E17").Value = "-2.46%"
$ws.Range("G17").Value = "4"
$ws.Range("E18").Value = "2.46%"
$ws.Range("G18").Value = "4"
$ws.Range("D19").Value = "0.3485"
$ws.Range("E19").Value = "-0.70%"
$ws.Range("G19").Value = "4"
$ws.Range("D20").Value = "8.460"
$ws.Range("E20").Value = "-7.47%"
$ws.Range("G20").Value = "4"
$ws.Range("D21").Value = "0.1382"
$ws.Range("E21").Value = "0.74%"
$ws.Range("G21").Value = "4"
$ws.Range("E22").Value = "3.45%"
$ws.Range("G22").Value = "4"
$ws.Range("D23").Value = "0.04271"
$ws.Range("E23").Value = "-2.65%"
$ws.Range("G23").Value = "4"
$ws.Range("E24").Value = "-1.61%"
$ws.Range("G24").Value = "4"
$ws.Range("D25").Value = "0.004650"
$ws.Range("E25").Value = "4.99%"
$ws.Range("G25").Value = "4"
$ws.Range("D26").Value = "0.0001233"
$ws.Range("E26").Value = "3.42%"
$ws.Range("G26").Value = "4"
$ws.Range("D27").Value = "0.0003999"
$ws.Range("E27").Value = "0.05%"
$ws.Range("G27").Value = "4"
$ws.Range("G28").Value = "4"
$ws.Range("G29").Value = "4"
$ws.Range("G30").Value = "4"
$ws.Range("G31").Value = "4"
$ws.Range("G32").Value = "4"
$ws.Range("G33").Value = "4"
$ws.Range("G34").Value = "4"
$ws.Range("G35").Value = "4"
$ws.Range("G36").Value = "4"
$ws.Range("G37").Value = "4"
$ws.Range("G38").Value = "4"
$ws.Range("D39").Value = "0.02662"
$ws.Range("E39").Value = "-6.13%"
$ws.Range("G39").Value = "4"
$ws.Range("D40").Value = "0.05556"
$ws.Range("E40").Value = "-1.01%"
$ws.Range("G40").Value = "4"
$ws.Range("E41").Value = "24.75%"
$ws.Range("G41").Value = "4"
$ws.Range("D42").Value = "0.007703"
$ws.Range("E42").Value = "-2.79%"
$ws.Range("G42").Value = "4"
$ws.Range("D43").Value = "0.1393"
$ws.Range("G43").Value = "4"
$ws.Range("D44").Value = "0.002059"
$ws.Range("E44").Value = "-2.15%"
$ws.Range("G44").Value = "4"
$ws.Range("D45").Value = "0.008695"
$ws.Range("E45").Value = "-11.85%"
$ws.Range("G45").Value = "4"
$ws.Range("D46").Value = "0.00007111"
$ws.Range("E46").Value = "-2.94%"
$ws.Range("G46").Value = "4"
$ws.Range("E47").Value = "0.05%"
$ws.Range("G47").Value = "4"
$ws.Range("D48").Value = "0.003494"
$ws.Range("E48").Value = "-4.93%"
$ws.Range("G48").Value = "4"
$ws.Range("D49").Value = "0.002276"
$ws.Range("E49").Value = "-0.28%"
$ws.Range("G49").Value = "4"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("E50").Value = "0.05%"
$ws.Range("G50").Value = "4"
$ws.Range("D51").Value = "0.0002005"
$ws.Range("E51").Value = "0.05%"
$ws.Range("G51").Value = "4"

# Restore the default (unstyled) cell style so no extra formatting is introduced
$changeRange.Style = "Normal"
